# Add 8 new workout rows (rows 244-251) for the week of 7/20/2024,
# matching the "Add files via upload" commit that appended fresh
# Strava export data to the Kilimanjaro weekly scoreboard sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (Participant, Date serial, Workout Type, Total Duration,
# Total Distance, Total Elevation, Zone1..Zone5, Workout Level, Week)
$newRows = @(
    @{ A="Steven"; B=45493; C="Walk";    D=23;   E=1.1499999999999999; F=56;  G=23; H=0;  I=0;  J=0; K=0; L="Brave Leopard";   M=6 },
    @{ A="Steven"; B=45493; C="Walk";    D=27;   E=1.2;                F=59;  G=27; H=0;  I=0;  J=0; K=0; L="Brave Leopard";   M=6 },
    @{ A="Matt";   B=45493; C="Run";     D=74;   E=6.01;               F=561; G=3;  H=35; I=20; J=6; K=0; L="Wily Hyena";      M=6 },
    @{ A="Matt";   B=45493; C="Walk";    D=0.25; E=6;                  F=0;   G=6;  H=0;  I=0;  J=0; K=0; L="Wily Hyena";      M=6 },
    @{ A="Phil";   B=45493; C="Workout"; D=61;   E=0;                  F=0;   G=7;  H=38; I=16; J=0; K=0; L="Sauntering Hippo"; M=6 },
    @{ A="Steven"; B=45493; C="Walk";    D=44;   E=2.27;               F=89;  G=36; H=8;  I=0;  J=0; K=0; L="Brave Leopard";   M=6 },
    @{ A="Eric";   B=45493; C="Workout"; D=45;   E=0;                  F=0;   G=45; H=0;  I=0;  J=0; K=0; L="Wily Hyena";      M=6 },
    @{ A="Steven"; B=45493; C="Walk";    D=28;   E=1.38;               F=56;  G=27; H=1;  I=0;  J=0; K=0; L="Brave Leopard";   M=6 }
)

$startRow = 244
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    # Seed the row by copying the last existing data row (243), which
    # already carries the date-formatted style for column B, so the new
    # rows reuse the same cell style instead of minting a new one.
    $ws.Range("A243:M243").Copy($ws.Range("A" + $r + ":M" + $r))

    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
}

# Keep the active selection pointed at the new bottom of the table,
# matching how Excel nudges the view after appending rows.
$ws.Range("A252").Select() | Out-Null
